# Auto-generated: update TOTAL AIRCRAFT CURVES sheet data (CL total / CM total rows)
# Commit: More work on ACAerodynamicAndStabilityManager.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOTAL AIRCRAFT CURVES")

$startCol = 3  # column C

$rowValues = @(0.4685493483143862, 0.5428473841585919, 0.6171454200027975, 0.6914434558470032, 0.7657414916912088, 0.8400395275354144, 0.9143375633796201, 0.9886355992238257, 1.0629336350680314, 1.137231670912237, 1.2115297067564428, 1.2858277426006484, 1.360125778444854, 1.4344238142890595, 1.5087218501332653, 1.5830198859774711, 1.6573179218216767, 1.7316160175582527, 1.8059695711958836, 1.880555788594429, 1.9554894996834373, 2.0307439340028313, 2.1047016416008977, 2.1741792325004043, 2.2364572570184187, 2.2894072863012522, 2.331491912324421, 2.3617647478927095, 2.3798704266401276, 2.386044603029913, 2.381113952354558, 2.3664961707357963, 2.3441999751245812, 2.3168251033011087)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(6, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(-0.22913350696062929, -0.2913803296109893, -0.35360197038052843, -0.4157847153663583, -0.47784368057090343, -0.5394173768526624, -0.5827565579596732, -0.6257172710832709, -0.6686918330151035, -0.7116904406772329, -0.7547215127791411, -0.7977945927563481, -0.8409161978125499, -0.8840963639986803, -0.9273448300395607, -0.9706683274264529, -1.0140733804497046, -1.0575706252378834, -1.101512653443434, -1.147004327011349, -1.1947801422649658, -1.2460586328206635, -1.312009105237918, -1.4067488158572634, -1.5013770482389646, -1.5666595858263772, -1.6250523531184031, -1.6736549513871093, -1.7086425010724207, -1.724834201551132, -1.7157178582843184, -1.6766577007446077, -1.6016353483910881, -1.4840329427216785)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(10, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(-0.18382925496526886, -0.23845080036011113, -0.29303498666588246, -0.34759051964352533, -0.40203270816640646, -0.4560025654632626, -0.49175501004771327, -0.5271465346955796, -0.5625711374454146, -0.5980414859379146, -0.6335685421034954, -0.6691631037234078, -0.7048341679087878, -0.7405938260023681, -0.7764526275534311, -0.8124206353897175, -0.8485068222263294, -0.8847224455894883, -0.9214158102047108, -0.9596782935021855, -1.0002353501850734, -1.0443066600749473, -1.1032228871957472, -1.1914149055939014, -1.280247155582478, -1.3407825256823087, -1.3956489693852019, -1.4419995569734116, -1.4759985612018982, -1.492426710356614, -1.4846768041391885, -1.4479731882778648, -1.3761088333490137, -1.2622293420777901)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(13, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(-0.147585853368981, -0.19610717695940877, -0.24458139969416628, -0.2930351630652597, -0.3413839302428096, -0.38927071635174354, -0.4189537717181461, -0.44828994558542756, -0.47767458098966464, -0.507122322146461, -0.5366461655629803, -0.5662579124970565, -0.5959685439857793, -0.6257917956053196, -0.6557388655645291, -0.6858224817603304, -0.7160535756476307, -0.7464439018707739, -0.7773383356137334, -0.8098174666948565, -0.8445995165211613, -0.8829050818783761, -0.9361939127620126, -1.0191477773832132, -1.1033432414572903, -1.1600808775670555, -1.2121262623986424, -1.2566752414424551, -1.2898834093054825, -1.3065007174010015, -1.299843960823086, -1.2650255783044726, -1.195687621315356, -1.0847864615626817)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(16, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.43525195391506993, 0.5103344507214107, 0.5854169475277512, 0.6604994443340919, 0.7355819411404325, 0.8106644379467732, 0.885746934753114, 0.9608294315594545, 1.0359119283657952, 1.1109944251721358, 1.1860769219784766, 1.2611594187848172, 1.3362419155911578, 1.4113244123974984, 1.4864069092038392, 1.56148940601018, 1.6365719028165207, 1.7116543996228613, 1.7867376037269689, 1.861908840056377, 1.9373735262671827, 2.013196503661717, 2.0878277237145557, 2.1581269875386444, 2.2213938262051123, 2.275494581278217, 2.318862404815306, 2.3504972593668843, 2.369965917976562, 2.3774019641810713, 2.373505792010276, 2.359544605987165, 2.3373524211278403, 2.30933006294153)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(20, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(-0.020733305440624572, -0.08784817210832285, -0.15493785689520045, -0.22198864589836886, -0.2889156551202524, -0.35535739541934985, -0.40356462054369896, -0.45139337768463506, -0.49923598363380606, -0.5471026353132737, -0.5950017514325203, -0.6429428754270661, -0.6909325245006062, -0.738980734704075, -0.787097244762294, -0.8352887861665244, -0.8835618832071145, -0.931926800345067, -0.9803967410402503, -1.0295188029179487, -1.0805896326971929, -1.1353962691379382, -1.2055262844952677, -1.305364961357414, -1.4061293700843795, -1.4785528414713276, -1.5449086079210195, -1.6019633391614687, -1.6454089846405688, -1.6694313332952062, -1.6667345406985508, -1.6317489469873807, -1.5573720428370468, -1.4357516015137761)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(24, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.02124120711480415, -0.03816993620116277, -0.0975437204280592, -0.15688885132682706, -0.21612063777083304, -0.2748800929888141, -0.31542213549438963, -0.355603258063381, -0.3958174587343408, -0.43607740514796567, -0.47639405923467154, -0.5167782187757088, -0.5572388808822137, -0.5977881368969189, -0.6384365363691074, -0.6791941421265182, -0.7200699268842548, -0.7610747824902113, -0.8022230945484183, -0.8440574642625904, -0.8878564379589258, -0.9353990394263336, -0.9984274582417312, -1.0916362755902278, -1.1865058205092234, -1.2540670518295625, -1.3167681749387292, -1.3714346936003534, -1.413755495636403, -1.4378881059855726, -1.436454302587849, -1.4037595909955014, -1.3325302831946466, -1.214697504905846)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(27, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.05482081715914677, 0.0015726525245650658, -0.051628411254346696, -0.10480901566959422, -0.15788462389129831, -0.2104982510443863, -0.244908147454943, -0.2789711623663786, -0.31308263881476966, -0.3472572210157202, -0.3815079054763935, -0.41584649345462404, -0.450283965987501, -0.4848340586511952, -0.5195079696545591, -0.5543184268945145, -0.5892763618259687, -0.6243931682063282, -0.6596841773549541, -0.6956883933383056, -0.7336698821683141, -0.7754012556570515, -0.8327483972389035, -0.9206533269764805, -1.0108069808491003, -1.0744784201161524, -1.1342558285528992, -1.187011777151463, -1.2284327044330723, -1.2526535241378673, -1.2522301120992887, -1.2213681062019994, -1.1526568754807278, -1.0378542276195035)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(30, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.40775208263693175, 0.4833615421727999, 0.558971001708668, 0.6345804612445363, 0.7101899207804043, 0.7857993803162725, 0.8614088398521407, 0.9370182993880087, 1.012627758923877, 1.088237218459745, 1.1638466779956134, 1.2394561375314814, 1.3150655970673495, 1.3906750566032176, 1.466284516139086, 1.5418939756749543, 1.6175034352108224, 1.6931128947466905, 1.7687223542825588, 1.844353937698859, 1.920199275032274, 1.9963953475505891, 2.0714476719012, 2.1422539608786524, 2.206134003744465, 2.2609567467620573, 2.3051402931967244, 2.3376519033157086, 2.358007994388132, 2.366274140685018, 2.3630650734793113, 2.349544681045862, 2.3274260086614103, 2.2989712586046203)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(34, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.15235693939332196, 0.08197195746602526, 0.011612157419549796, -0.0587087468432167, -0.12890587132469822, -0.1986177268833938, -0.2500950672673409, -0.301193939667875, -0.352306660876644, -0.40344342781570985, -0.4546126591945546, -0.5058238984486985, -0.5570836627818367, -0.6084019882449038, -0.6597886135627207, -0.7112502702265493, -0.7627934825267376, -0.814428514924288, -0.866164181678295, -0.9181429693720804, -0.9715759643954287, -1.02869787664069, -1.1014410931085854, -1.204426160783799, -1.3089958640736752, -1.3858996974714848, -1.4573175067289597, -1.519813024573962, -1.5687657131050217, -1.5979393235604846, -1.5995064234658607, -1.5672562145610844, -1.493335498775474, -1.3690318102801535)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(38, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.19158146482093685, 0.12895290251832425, 0.06636169930478275, 0.0037991494193695408, -0.058650056011281704, -0.12062693021590815, -0.16438639170812894, -0.20778493326376546, -0.25121655292137063, -0.2946939183216408, -0.33822799139499216, -0.3818295699226748, -0.4255076510158251, -0.4692743260171756, -0.5131401444760093, -0.5571151692200657, -0.601208372964448, -0.6454306475570494, -0.6897920601309039, -0.7344371209524743, -0.7805601947806529, -0.8303807625401984, -0.8959802720363843, -0.992284777682612, -1.0908982967445837, -1.1628676912813356, -1.2305492849085278, -1.2905689146179642, -1.3383080164596994, -1.3675088786004561, -1.3702702572082552, -1.340266851063335, -1.2694863803797167, -1.149013594105914)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(41, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.22296108516302826, 0.16653765856016317, 0.11016133281296858, 0.0538054664294379, -0.0024454037605492474, -0.05823429288192039, -0.09581945126076011, -0.13305772814047873, -0.17034446655715296, -0.2076943107263866, -0.24512025715534314, -0.28263410710185694, -0.320246841603017, -0.35797219623499454, -0.39582136920664157, -0.4338070884148801, -0.4719402853146175, -0.5102323536632598, -0.5486943628929927, -0.5874724422167907, -0.6277475790888336, -0.6717270712598065, -0.7316116151786255, -0.8225716712016642, -0.9164202428813126, -0.9844420863292183, -1.049134707452184, -1.107173626653168, -1.153941859143443, -1.1831645226324352, -1.1868813242021723, -1.1586753602651376, -1.090407085663112, -0.9729990211665245)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(44, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.40322933381319764, 0.4785995099634247, 0.5539696861136516, 0.6293398622638787, 0.7047100384141057, 0.7800802145643327, 0.8554503907145599, 0.9308205668647869, 1.0061907430150139, 1.0815609191652409, 1.1569310953154681, 1.2323012714656951, 1.307671447615922, 1.383041623766149, 1.4584117999163762, 1.5337819760666034, 1.6091521522168304, 1.6845223283670572, 1.759892513464941, 1.8353214348550053, 1.9110212166074383, 1.987082528182125, 2.061974129518935, 2.1325676105038465, 2.196167271288762, 2.250637202826447, 2.2944012868705017, 2.3264431959754144, 2.3463063934965263, 2.3540941335900296, 2.3504694612129935, 2.3366552121233495, 2.314434012879882, 2.286148280842234)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(48, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.18342937480336374, 0.11452928774036808, 0.04565438255819332, -0.023181626840272473, -0.09189385645745324, -0.1601208171518479, -0.2101132626714942, -0.2597272402077275, -0.3093550665521958, -0.35900693862696076, -0.4086912751415047, -0.45841761953134763, -0.5081924890001851, -0.5580259195989513, -0.6079276500524674, -0.657904411851995, -0.7079627292878827, -0.7581128668211324, -0.8083636942363366, -0.8590848442502648, -0.9116145808459875, -0.9679002217698374, -1.0396460569229822, -1.1413105251994575, -1.244140291924785, -1.3188547215193318, -1.3876695186708075, -1.4472502640518417, -1.493144256714706, -1.5193490738603037, -1.5183370979855937, -1.4842633348108265, -1.4097063792545448, -1.286451547485846)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(52, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.22220162534860521, 0.16103402957172952, 0.09990379288392459, 0.038802209524247955, -0.02218602938066655, -0.0827019370595562, -0.1250004320260403, -0.16693800705594025, -0.20890866018780874, -0.25092505906234214, -0.2929981656099568, -0.3351387776119026, -0.3773558921793163, -0.41966160065493, -0.4620664525880269, -0.5045805108063465, -0.5472127480249922, -0.589974056091857, -0.6328745567707073, -0.6762822461150437, -0.7215166170736951, -0.7705143896061919, -0.8351325900890079, -0.9301377771357509, -1.0270393978412644, -1.0968546697227435, -1.161975197482998, -1.2191270248298733, -1.2638567201585436, -1.290136629609774, -1.2903604929546204, -1.2585629182053282, -1.1871564604369405, -1.067715629087845)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(55, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.25321942578479806, 0.19823782303681842, 0.14330332114450917, 0.08838927861586388, 0.03358023228076207, -0.020766832985723627, -0.056910167509678006, -0.09270662053451126, -0.1285515350963001, -0.16445955541064824, -0.20044367798471946, -0.2365157040763477, -0.2726866147226224, -0.3089701454997144, -0.34537749461647593, -0.3819213899698289, -0.4186127630146812, -0.45546300750843816, -0.49248324679820543, -0.5300401676068685, -0.5694382460558629, -0.6126057238752773, -0.6715218166218303, -0.7611995786847874, -0.8533586825744495, -0.919254628285475, -0.9814197405327518, -1.0366284334523004, -1.080426690913616, -1.1067666742093523, -1.107979208929843, -1.078002584920932, -1.0091165253828587, -0.8927268943694464)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(58, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.396683776367474, 0.4718389560105058, 0.5469941356535375, 0.6221493152965694, 0.697304494939601, 0.7724596745826329, 0.8476148542256647, 0.9227700338686965, 0.9979252135117282, 1.0730803931547601, 1.148235572797792, 1.2233907524408236, 1.2985459320838555, 1.3737011117268871, 1.448856291369919, 1.524011471012951, 1.5991666506559827, 1.6743218302990146, 1.7494797642034439, 1.8247522663206335, 1.9003378507706106, 1.976285252963975, 2.0510309906668724, 2.1214241034671164, 2.1847560330940037, 2.238887703953261, 2.2822495231270024, 2.313841380373798, 2.3332326481286265, 2.3405621815028814, 2.336538318284391, 2.322438878937414, 2.3001111666026155, 2.2719719670970795)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(62, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.2274377979085419, 0.1598718912714936, 0.09233116651526618, 0.024829337542748076, -0.042548711648485116, -0.10944149191693242, -0.15809975701063123, -0.20637955412091705, -0.25467320003943783, -0.30299089168825544, -0.3513410477768519, -0.39973321174074733, -0.4481739007836373, -0.49667315095645603, -0.5452407009840248, -0.5938832823576048, -0.6426074193675452, -0.6914233764748474, -0.7403570597588052, -0.7901075357146605, -0.8419286105389894, -0.8975073764122876, -0.9683480414872663, -1.0687691067268563, -1.1699374433428429, -1.2425527654592563, -1.3088712712847936, -1.3656591788304868, -1.4086245683043446, -1.4319859418803982, -1.4284967508460844, -1.3926532135295502, -1.3174352812745553, -1.195089770492779)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(66, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.26555549270921086, 0.20570057770756311, 0.14588302179498613, 0.08609411921053767, 0.026418561080851055, -0.0327846658228107, -0.07377048001406683, -0.1143953742687388, -0.15505334662537934, -0.19575706472468493, -0.23651749049707155, -0.27734542172378945, -0.31824985551597507, -0.35924288321636105, -0.4003350543742301, -0.4415364318173216, -0.48285598826073944, -0.5243046155523763, -0.5659091972193254, -0.6083618544328766, -0.65289898335038, -0.7012012717704572, -0.7649288885384982, -0.8587107093668229, -0.9539776730787979, -1.0217276635499866, -1.084392126471334, -1.1387961211686803, -1.1806444062849721, -1.2041266928385834, -1.201913260107971, -1.168374430242646, -1.0963176470846774, -0.9777714834692939)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(69, $startCol + $i).Value = $rowValues[$i]
}

$rowValues = @(0.2960496485497457, 0.24236352685641854, 0.18872450601876156, 0.13510594454476865, 0.08159237926431928, 0.028540795052485857, -0.006307058416816169, -0.04080803038699707, -0.07535746389413354, -0.1099700031538295, -0.14465864467324827, -0.1794351897102242, -0.2143106193018465, -0.24929866902428624, -0.28441053708639563, -0.31965895138509626, -0.35505484337529625, -0.3906096068144008, -0.4263509071877434, -0.46296530940745123, -0.5016752815994939, -0.5441563880569946, -0.6021935661794854, -0.690663991478798, -0.7812098568675636, -0.8450675820225728, -0.9048088106205683, -0.957305675039237, -0.9982602766694764, -1.0218392936051335, -1.020646467517482, -0.9889514036131242, -0.9194235397327769, -0.8039168538505073)
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(72, $startCol + $i).Value = $rowValues[$i]
}

